$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2979
$ws.Range("I19").Value = 4986
$ws.Range("J19").Value = 1473.75
$ws.Range("K19").Value = 4986
$ws.Range("L19").Value = 1473.75
$ws.Range("M19").Value = -4811
$ws.Range("N19").Value = -1823.75
$ws.Range("H28").Value = 3604.5
$ws.Range("I28").Value = 632.875
$ws.Range("J28").Value = 9547.75
$ws.Range("K28").Value = 632.875
$ws.Range("L28").Value = 9547.75
$ws.Range("M28").Value = -147.875
$ws.Range("N28").Value = -10517.75
$ws.Range("H29").Value = 3218.182
$ws.Range("J29").Value = 5180
$ws.Range("L29").Value = 15540
$ws.Range("N29").Value = -16102
$ws.Range("H40").Value = 8179.089
$ws.Range("I40").Value = 6530.2666
$ws.Range("J40").Value = 9003.5
$ws.Range("K40").Value = 6530.2666
$ws.Range("L40").Value = 9003.5
$ws.Range("M40").Value = -6355.2666
$ws.Range("N40").Value = -9353.5
$ws.Range("H64").Value = 9416.799999999999
$ws.Range("J64").Value = 12867.333
$ws.Range("L64").Value = 12867.333
$ws.Range("N64").Value = -13363.333
$ws.Range("H67").Value = 9416.799999999999
$ws.Range("J67").Value = 12867.333
$ws.Range("L67").Value = 12867.333
$ws.Range("N67").Value = -14583.333
$ws.Range("H74").Value = 9600.733
$ws.Range("I74").Value = 8650.700000000001
$ws.Range("K74").Value = 8650.700000000001
$ws.Range("M74").Value = -7714.700000000001
$ws.Range("H77").Value = 9600.733
$ws.Range("I77").Value = 8650.700000000001
$ws.Range("K77").Value = 43253.5
$ws.Range("M77").Value = -38573.5
$ws.Range("H132").Value = 3434.8823
$ws.Range("I132").Value = 2026.8966
$ws.Range("J132").Value = 11601.2
$ws.Range("K132").Value = 6080.6898
$ws.Range("L132").Value = 34803.60000000001
$ws.Range("M132").Value = -3550.6898
$ws.Range("N132").Value = -39863.60000000001
$ws.Range("H135").Value = 1810.9286
$ws.Range("I135").Value = 1720
$ws.Range("J135").Value = 2993
$ws.Range("K135").Value = 15480
$ws.Range("L135").Value = 26937
$ws.Range("M135").Value = -12945
$ws.Range("N135").Value = -32007
$ws.Range("H137").Value = 3290.2122
$ws.Range("I137").Value = 2214.95
$ws.Range("J137").Value = 4944.4614
$ws.Range("K137").Value = 6644.849999999999
$ws.Range("L137").Value = 14833.3842
$ws.Range("M137").Value = -4094.849999999999
$ws.Range("N137").Value = -19933.3842
$ws.Range("H138").Value = 5004.243
$ws.Range("J138").Value = 4710.909
$ws.Range("L138").Value = 14132.727
$ws.Range("N138").Value = -24412.727
$ws.Range("H140").Value = 70433.8
$ws.Range("J140").Value = 70433.8
$ws.Range("L140").Value = 70433.8
$ws.Range("N140").Value = -80793.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5999.75
$ws.Range("I6").Value = 5999
$ws.Range("K6").Value = 5999
$ws.Range("M6").Value = -5826

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2099.75
$ws.Range("I7").Value = 199.5
$ws.Range("K7").Value = 199.5
$ws.Range("M7").Value = -86.5
$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""
$ws.Range("H31").Value = 6750
$ws.Range("J31").Value = 9875
$ws.Range("L31").Value = 9875
$ws.Range("N31").Value = -10379
$ws.Range("H94").Value = 1351.2
$ws.Range("I94").Value = 1351.2
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1351.2
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -900.2
$ws.Range("H107").Value = 1321.2727
$ws.Range("I107").Value = 1314.8889
$ws.Range("K107").Value = 1314.8889
$ws.Range("M107").Value = 605.1111000000001
$ws.Range("H134").Value = 4111.353
$ws.Range("I134").Value = 2353.889
$ws.Range("K134").Value = 7061.667
$ws.Range("M134").Value = -4526.667
$ws.Range("H140").Value = 66236.22
$ws.Range("J140").Value = 66236.22
$ws.Range("L140").Value = 66236.22
$ws.Range("N140").Value = -76596.22

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 3700
$ws.Range("J33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6758
$ws.Range("H38").Value = 15956.714
$ws.Range("I38").Value = 9174.5
$ws.Range("J38").Value = 24999.666
$ws.Range("K38").Value = 9174.5
$ws.Range("L38").Value = 24999.666
$ws.Range("M38").Value = -8797.5
$ws.Range("N38").Value = -25753.666
$ws.Range("H46").Value = 15956.714
$ws.Range("I46").Value = 9174.5
$ws.Range("J46").Value = 24999.666
$ws.Range("K46").Value = 9174.5
$ws.Range("L46").Value = 24999.666
$ws.Range("M46").Value = -8963.5
$ws.Range("N46").Value = -25421.666
$ws.Range("H134").Value = 4586.8945
$ws.Range("I134").Value = 2247.9092
$ws.Range("K134").Value = 6743.7276
$ws.Range("M134").Value = -4208.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 495.2
$ws.Range("J23").Value = 594
$ws.Range("L23").Value = 1782
$ws.Range("N23").Value = -2252
$ws.Range("H106").Value = 18014.5
$ws.Range("J106").Value = 18014.5
$ws.Range("L106").Value = 54043.5
$ws.Range("N106").Value = -55935.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 849.75
$ws.Range("I41").Value = 849.75
$ws.Range("K41").Value = 849.75
$ws.Range("M41").Value = -494.75
$ws.Range("H53").Value = 36642.145
$ws.Range("J53").Value = 29247.5
$ws.Range("L53").Value = 29247.5
$ws.Range("N53").Value = -30509.5
$ws.Range("H54").Value = 4095
$ws.Range("J54").Value = 4095
$ws.Range("L54").Value = 4095
$ws.Range("N54").Value = -4875
$ws.Range("H97").Value = 2445.3333
$ws.Range("I97").Value = 1560.5555
$ws.Range("K97").Value = 1560.5555
$ws.Range("M97").Value = -1064.5555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 7000
$ws.Range("J33").Value = 7000
$ws.Range("L33").Value = 7000
$ws.Range("N33").Value = -7580
$ws.Range("H53").Value = 35868.625
$ws.Range("J53").Value = 32089.8
$ws.Range("L53").Value = 32089.8
$ws.Range("N53").Value = -33125.8
$ws.Range("H61").Value = 2578.2856
$ws.Range("I61").Value = 2186.3635
$ws.Range("K61").Value = 2186.3635
$ws.Range("M61").Value = -1984.3635
$ws.Range("H101").Value = 11117.333
$ws.Range("J101").Value = 11117.333
$ws.Range("L101").Value = 11117.333
$ws.Range("N101").Value = -17607.333
$ws.Range("H104").Value = 15185
$ws.Range("J104").Value = 15185
$ws.Range("L104").Value = 15185
$ws.Range("N104").Value = -22173
$ws.Range("H113").Value = 2578.2856
$ws.Range("I113").Value = 2186.3635
$ws.Range("K113").Value = 2186.3635
$ws.Range("M113").Value = -16.36349999999993
$ws.Range("H132").Value = 4504.778
$ws.Range("J132").Value = 7361.25
$ws.Range("L132").Value = 22083.75
$ws.Range("N132").Value = -27143.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 6663.3335
$ws.Range("J24").Value = 7495
$ws.Range("L24").Value = 7495
$ws.Range("N24").Value = -7955
$ws.Range("H41").Value = 18616
$ws.Range("J41").Value = 19890
$ws.Range("L41").Value = 19890
$ws.Range("N41").Value = -20670
$ws.Range("H58").Value = 7816.6665
$ws.Range("I58").Value = 6000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5692
$ws.Range("H103").Value = 37812.5
$ws.Range("J103").Value = 37812.5
$ws.Range("L103").Value = 37812.5
$ws.Range("N103").Value = -40156.5
$ws.Range("H104").Value = 14210
$ws.Range("J104").Value = 14210
$ws.Range("L104").Value = 14210
$ws.Range("N104").Value = -21198
$ws.Range("H113").Value = 1590.9667
$ws.Range("J113").Value = 1763.2222
$ws.Range("L113").Value = 5289.6666
$ws.Range("N113").Value = -9629.6666
$ws.Range("H132").Value = 5138.766
$ws.Range("I132").Value = 3450.6316
$ws.Range("K132").Value = 10351.8948
$ws.Range("M132").Value = -7821.8948

Write-Output "Applied 207 cell changes across 8 sheets"